# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.648.40"
$ws.Range("E2").Value = "  +5.14%  "

$ws.Range("D3").Value = "'1.920.96"
$ws.Range("E3").Value = "  +3.74%  "

$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").Value = "'335.23"
$ws.Range("E5").Value = "  +1.37%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.19%  "

$ws.Range("D7").Value = "'0.4679"
$ws.Range("E7").Value = "  +3.14%  "

$ws.Range("D8").Value = "'0.4116"
$ws.Range("E8").Value = "  +5.34%  "

$ws.Range("D9").Value = "'48.20"
$ws.Range("E9").Value = "  +0.93%  "

$ws.Range("D10").Value = "'0.08037"
$ws.Range("E10").Value = "  +3.06%  "

$ws.Range("D11").Value = "'1.015"
$ws.Range("E11").Value = "  +3.91%  "

$ws.Range("D12").Value = "'22.45"
$ws.Range("E12").Value = "  +5.07%  "

$ws.Range("D13").Value = "'1.944.50"
$ws.Range("E13").Value = "  +5.42%  "

$ws.Range("D14").Value = "'5.999"
$ws.Range("E14").Value = "  +3.64%  "

$ws.Range("D15").Value = "'7.187"
$ws.Range("E15").Value = "  +3.34%  "

$ws.Range("D16").Value = "'89.94"
$ws.Range("E16").Value = "  +2.87%  "

$ws.Range("E17").Value = "  -0.12%  "

$ws.Range("D18").Value = "'0.00001037"
$ws.Range("E18").Value = "  +2.11%  "

$ws.Range("D19").Value = "'0.06593"
$ws.Range("E19").Value = "  +0.95%  "

$ws.Range("D20").Value = "'17.83"
$ws.Range("E20").Value = "  +5.17%  "

$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  -1.48%  "

$ws.Range("D22").Value = "'29.609.30"
$ws.Range("E22").Value = "  +5.07%  "

$ws.Range("D23").Value = "'5.578"
$ws.Range("E23").Value = "  +5.93%  "

$ws.Range("D24").Value = "'11.69"
$ws.Range("E24").Value = "  +10.46%  "

$ws.Range("D25").Value = "'2.207"
$ws.Range("E25").Value = "  -2.10%  "

$ws.Range("D26").Value = "'2.170.32"
$ws.Range("E26").Value = "  +5.16%  "

$ws.Range("D27").Value = "'156.03"
$ws.Range("E27").Value = "  -0.04%  "

$ws.Range("D28").Value = "'19.87"
$ws.Range("E28").Value = "  +4.01%  "

$ws.Range("E29").Value = "  +5.74%  "

$ws.Range("D30").Value = "'5.710"
$ws.Range("E30").Value = "  +8.74%  "

$ws.Range("D31").Value = "'117.63"
$ws.Range("E31").Value = "  +1.23%  "

$ws.Range("D32").Value = "'1.076"
$ws.Range("E32").Value = "  +15.33%  "

$ws.Range("D33").Value = "'0.09476"
$ws.Range("E33").Value = "  +2.55%  "

$ws.Range("E34").Value = "  +4.47%  "

$ws.Range("D35").Value = "'3.576"
$ws.Range("E35").Value = "  -0.82%  "

$ws.Range("D36").Value = "'5.420"
$ws.Range("E36").Value = "  +4.80%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02276"
$ws.Range("E37").Value = "  +4.19%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.06136"
$ws.Range("E38").Value = "  +2.11%  "

$ws.Range("D39").Value = "'8.439"
$ws.Range("E39").Value = "  +3.62%  "

$ws.Range("E40").Value = "  +1.26%  "

$ws.Range("D41").Value = "'0.5893"
$ws.Range("E41").Value = "  +4.61%  "

$ws.Range("D42").Value = "'0.1846"
$ws.Range("E42").Value = "  +3.53%  "

$ws.Range("D43").Value = "'10.22"
$ws.Range("E43").Value = "  +2.64%  "

$ws.Range("D44").Value = "'1.259"
$ws.Range("E44").Value = "  +1.31%  "

$ws.Range("D45").Value = "'2.344"
$ws.Range("E45").Value = "  +1.96%  "

$ws.Range("D46").Value = "'0.07504"
$ws.Range("E46").Value = "  +4.43%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'12.27"
$ws.Range("E47").Value = "  +5.01%  "

$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "'0.5584"
$ws.Range("E48").Value = "  +4.29%  "

$ws.Range("D49").Value = "'1.937"
$ws.Range("E49").Value = "  +3.88%  "

$ws.Range("D50").Value = "'113.16"
$ws.Range("E50").Value = "  +3.47%  "

$ws.Range("D51").Value = "'0.2998"
$ws.Range("E51").Value = "  +14.15%  "
